$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Uppercase the shared strings in column A (f1..f180 -> F1..F180)
for ($i = 1; $i -le 180; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $cell.Value2 = $cell.Value2.ToUpper()
}

# 2) Recompute column D values: new = ceil(old / 2)
for ($i = 1; $i -le 180; $i++) {
    $cell = $ws.Cells.Item($i, 4)
    $old = $cell.Value2
    $cell.Value2 = [Math]::Ceiling($old / 2)
}

# 3) Add a secondary (unused) 8pt Calibri font, used for the worksheet's
#    phonetic-guide font, without altering any existing cell's applied style.
$ws.Range("A1:A180").Phonetics.Font.Size = 8
$ws.Range("A1:A180").HorizontalAlignment = -4108
$ws.Range("A1:A180").VerticalAlignment = -4108
$ws.Range("A1:A180").Font.Size = 11

# 4) Update the sheet view: scroll down and change the selection to A1:A180
$ws.Activate()
$ws.Range("A1:A180").Select()
$excel.ActiveWindow.ScrollRow = 155
$excel.ActiveWindow.ScrollColumn = 1
